$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The title's text is already "Here is a single header" once concatenated
# across its (nine) runs, so a direct re-assignment of the identical string
# is treated as a no-op by the text engine. Round-trip through a throwaway
# value first so the engine actually rewrites the paragraph into a single
# run (collapsing "Here"/" "/"is"/... into one <a:r>).
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "__tmp__"
$titleRange.Text = "Here is a single header"

# Collapse the speaker notes' multi-run text ("and"/" "/"here"/...) into a
# single run as well.
$notes = $s.NotesPage
$notes.Shapes.Item(2).TextFrame.TextRange.Text = "and here are some notes"
